$wb = $excel.ActiveWorkbook

$configWs = $wb.Worksheets.Item("Config")
$testWs = $wb.Worksheets.Item("Test Cases")

# RunInParallel: Yes -> No
$configWs.Range("B3").Value = "No"

# TestCaseNumber filter: TestCaseNumber=All -> TestCaseNumber=101
$testWs.Range("D2").Value = "TestCaseNumber=101"

# Update selections to match final state
[void]$testWs.Activate()
[void]$testWs.Range("D2").Select()

[void]$configWs.Activate()
[void]$configWs.Range("C3").Select()

[void]$testWs.Activate()
